# Add a new, centered paragraph containing a YouTube link as the very
# first paragraph of the document body (before the existing cone-graphic
# image paragraph), matching:
#
#   <w:p>
#     <w:pPr>
#       <w:jc w:val="center"/>
#       <w:rPr>
#         <w:noProof/>
#       </w:rPr>
#     </w:pPr>
#     <w:r>
#       <w:rPr>
#         <w:noProof/>
#       </w:rPr>
#       <w:t>https://youtu.be/5Y_rzp0XOUs</w:t>
#     </w:r>
#   </w:p>

$d = $word.ActiveDocument

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/><w:rPr><w:noProof/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>https://youtu.be/5Y_rzp0XOUs</w:t></w:r></w:p>'

$insertionPoint = $d.Range(0, 0)
$insertionPoint.InsertXML($newParagraphXml)
